# Created experiment order generation script
# Regenerates the 5 task-order worksheets with freshly generated
# session identifiers, and reorders the "NB" / "RS" sheets.
#
# Final sheet order / names:
#   1. GNG_TO-16515889337923248
#   2. RS_TO-16515889337951927
#   3. NB_TO-16515889365943134
#   4. TOL_TO-1651588936642855
#   5. vSAT_TO-16515889367211888

$wb = $excel.ActiveWorkbook

# --- Step 1: reduce to a single worksheet (the GNG sheet), which already
#     carries the correct header / index-column styling (style id 1).
#     This sheet becomes the template that is copied for the other four
#     sheets so that every sheet ends up with matching formatting and a
#     clean, sequential sheetId (1..5) after the rebuild. ---
for ($i = $wb.Worksheets.Count; $i -ge 2; $i--) {
    [void]$wb.Worksheets.Item($i).Delete()
}

$template = $wb.Worksheets.Item(1)

# --- Step 2: create four more copies of the template sheet, appended
#     after the last sheet each time, so we end up with 5 identical
#     (styled) sheets in order. ---
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$excel.CutCopyMode = $false

$wsGNG  = $wb.Worksheets.Item(1)
$wsRS   = $wb.Worksheets.Item(2)
$wsNB   = $wb.Worksheets.Item(3)
$wsTOL  = $wb.Worksheets.Item(4)
$wsVSAT = $wb.Worksheets.Item(5)

# --- Step 3: rename sheets to their new names (new generated ids) ---
$wsGNG.Name  = "GNG_TO-16515889337923248"
$wsRS.Name   = "RS_TO-16515889337951927"
$wsNB.Name   = "NB_TO-16515889365943134"
$wsTOL.Name  = "TOL_TO-1651588936642855"
$wsVSAT.Name = "vSAT_TO-16515889367211888"

# --- GNG sheet: same shape as template (5 rows), only refresh the data
#     filenames. ---
$wsGNG.Range("B2").Value = "go_stims-16515889337614264.csv"
$wsGNG.Range("B3").Value = "GNG_stims-1651588933775271.csv"
$wsGNG.Range("B4").Value = "go_stims-16515889337772381.csv"
$wsGNG.Range("B5").Value = "GNG_stims-16515889337911866.csv"

# --- RS sheet: shrink from 5 rows down to 3 rows (header + 2 entries),
#     then set the (unchanged) "eyes closed" / "eyes open" values. ---
$wsRS.Range("A4:B5").Clear()
$wsRS.Range("B2").Value = "eyes closed"
$wsRS.Range("B3").Value = "eyes open"

# --- NB sheet: grow from 5 rows up to 10 rows (header + 9 entries).
#     Copy the index-column style from A2 down onto the new A6:A10
#     cells before filling in the values. ---
$wsNB.Range("A2").Copy()
$wsNB.Range("A6:A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsNB.Range("A6").Value = 4
$wsNB.Range("A7").Value = 5
$wsNB.Range("A8").Value = 6
$wsNB.Range("A9").Value = 7
$wsNB.Range("A10").Value = 8

$wsNB.Range("B2").Value  = "OB-16515889356107924.csv"
$wsNB.Range("B3").Value  = "TB-16515889357501245.csv"
$wsNB.Range("B4").Value  = "OB-16515889350314488.csv"
$wsNB.Range("B5").Value  = "OB-16515889346709244.csv"
$wsNB.Range("B6").Value  = "ZB-match_5-16515889340489855.csv"
$wsNB.Range("B7").Value  = "ZB-match_2-16515889343476653.csv"
$wsNB.Range("B8").Value  = "TB-16515889365810354.csv"
$wsNB.Range("B9").Value  = "ZB-match_6-1651588933995541.csv"
$wsNB.Range("B10").Value = "TB-1651588936038056.csv"

# --- TOL sheet: grow from 5 rows up to 7 rows (header + 6 entries).
#     Copy the index-column style from A2 down onto the new A6:A7
#     cells before filling in the values. ---
$wsTOL.Range("A2").Copy()
$wsTOL.Range("A6:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsTOL.Range("A6").Value = 4
$wsTOL.Range("A7").Value = 5

$wsTOL.Range("B2").Value = "MM_stims-165158893661007.csv"
$wsTOL.Range("B3").Value = "ZM_stims-1651588936596367.csv"
$wsTOL.Range("B4").Value = "MM_stims-16515889366267562.csv"
$wsTOL.Range("B5").Value = "ZM_stims-165158893661007.csv"
$wsTOL.Range("B6").Value = "MM_stims-16515889366418824.csv"
$wsTOL.Range("B7").Value = "ZM_stims-16515889366277637.csv"

# --- vSAT sheet: same shape as template (5 rows), only refresh the data
#     filenames. ---
$wsVSAT.Range("B2").Value = "SAT_stims-16515889366727197.csv"
$wsVSAT.Range("B3").Value = "vSAT_stims-16515889367053602.csv"
$wsVSAT.Range("B4").Value = "SAT_stims-1651588936644917.csv"
$wsVSAT.Range("B5").Value = "vSAT_stims-16515889366893275.csv"

$wsGNG.Select()
